$p = $ppt.ActivePresentation

# --- 1) Swap the table's style on slide 5 (the B1 financial-documents table)
#     from the deck-defined "Table_0" style to the built-in style PowerPoint
#     assigned when the table style was changed in the gallery.
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{15DDA280-BF0B-47F9-98F6-ACBE3C0B4E77}")

# --- 2) Switch the deck's design back to the default "Office Theme" colour
#     scheme (it was previously recoloured to the "Integral" / Red Violet
#     palette). Re-apply the 12 standard Office theme colours, in
#     msoThemeColorIndex order: Dark1, Light1, Dark2, Light2, Accent1-6,
#     Hyperlink, FollowedHyperlink.
$theme = $p.SlideMaster.Theme
$officeColors = @(
    0x000000,
    0xFFFFFF,
    0x44546A,
    0xE7E6E6,
    0x5B9BD5,
    0xED7D31,
    0xA5A5A5,
    0xFFC000,
    0x4472C4,
    0x70AD47,
    0x0563C1,
    0x954F72
)
for ($i = 1; $i -le $officeColors.Length; $i++) {
    $theme.ThemeColorScheme.Item($i).RGB = $officeColors[$i - 1]
}
